$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append two new days of data (rows 15 and 16), continuing the same
# --- date-serial / measurement layout used by the existing rows (2-14).
# --- Copy row 14's formatting down first so the new cells reuse the
# --- existing "date" (col A) and "bordered body" (col B:H) cell styles
# --- instead of minting brand-new ones.
$ws.Range("A14:H14").Copy() | Out-Null
$ws.Range("A15:H15").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:H16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(15, 1).Value = 45768
$ws.Cells.Item(15, 2).Value = 36.6
$ws.Cells.Item(15, 3).Value = 43.1
$ws.Cells.Item(15, 4).Value = 39.2
$ws.Cells.Item(15, 5).Value = 40.1
$ws.Cells.Item(15, 6).Value = 37.8
$ws.Cells.Item(15, 7).Value = 32.7
$ws.Cells.Item(15, 8).Value = 37.3

$ws.Cells.Item(16, 1).Value = 45769
$ws.Cells.Item(16, 2).Value = 36
$ws.Cells.Item(16, 3).Value = 41.6
$ws.Cells.Item(16, 4).Value = 38.7
$ws.Cells.Item(16, 5).Value = 40.5
$ws.Cells.Item(16, 6).Value = 37.5
$ws.Cells.Item(16, 7).Value = 31.6
$ws.Cells.Item(16, 8).Value = 37.9

# --- Highlight cell I14 with a right/top/bottom thin border (left edge
# --- left untouched/none), matching the new border style added to the
# --- workbook's style table.
$cell = $ws.Range("I14")
$cell.Borders.Item(7).LineStyle = 0    # xlEdgeLeft   -> no border
$cell.Borders.Item(10).LineStyle = 1   # xlEdgeRight  -> thin
$cell.Borders.Item(8).LineStyle = 1    # xlEdgeTop    -> thin
$cell.Borders.Item(9).LineStyle = 1    # xlEdgeBottom -> thin

# --- Leave the selection where the author ended up after the edit.
$ws.Range("D24").Select() | Out-Null
